$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($row, $col, $val) {
    $c = $ws.Cells.Item($row, $col)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

# Row 2
Set-TextValue 2 4 '55.998.29'
$ws.Cells.Item(2, 5).Value = '  -3.21%  '

# Row 3
Set-TextValue 3 4 '2.365.64'
$ws.Cells.Item(3, 5).Value = '  -2.38%  '

# Row 4
$ws.Cells.Item(4, 5).Value = '  -0.25%  '

# Row 5
Set-TextValue 5 4 '498.33'
$ws.Cells.Item(5, 5).Value = '  -2.34%  '

# Row 6
Set-TextValue 6 4 '130.33'
$ws.Cells.Item(6, 5).Value = '  -2.42%  '

# Row 7
Set-TextValue 7 4 '1.00'
$ws.Cells.Item(7, 5).Value = '  +0.19%  '

# Row 8
Set-TextValue 8 4 '0.545'
$ws.Cells.Item(8, 5).Value = '  -2.83%  '

# Row 9
Set-TextValue 9 4 '2.368.33'
$ws.Cells.Item(9, 5).Value = '  -3.67%  '

# Row 10
Set-TextValue 10 4 '0.0968'
$ws.Cells.Item(10, 5).Value = '  -1.85%  '

# Row 11
$ws.Cells.Item(11, 5).Value = '  +0.41%  '

# Row 12
Set-TextValue 12 4 '0.326'
$ws.Cells.Item(12, 5).Value = '  -0.12%  '

# Row 13
Set-TextValue 13 4 '4.63'
$ws.Cells.Item(13, 5).Value = '  -1.00%  '

# Row 14
Set-TextValue 14 4 '2.787.12'
$ws.Cells.Item(14, 5).Value = '  -2.57%  '

# Row 15
Set-TextValue 15 4 '55.934.19'
$ws.Cells.Item(15, 5).Value = '  -3.10%  '

# Row 16
Set-TextValue 16 4 '21.33'
$ws.Cells.Item(16, 5).Value = '  -3.34%  '

# Row 17
Set-TextValue 17 4 '0.0000131'
$ws.Cells.Item(17, 5).Value = '  -2.60%  '

# Row 18
Set-TextValue 18 4 '2.323.33'
$ws.Cells.Item(18, 5).Value = '  -5.60%  '

# Row 19
Set-TextValue 19 4 '9.99'
$ws.Cells.Item(19, 5).Value = '  -3.97%  '

# Row 20
$ws.Cells.Item(20, 2).Value = 'BitcoinCash'
$ws.Cells.Item(20, 3).Value = 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch'
Set-TextValue 20 4 '306.41'
$ws.Cells.Item(20, 5).Value = '  -3.33%  '

# Row 21
$ws.Cells.Item(21, 2).Value = 'Polkadot'
$ws.Cells.Item(21, 3).Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 21 4 '3.98'
$ws.Cells.Item(21, 5).Value = '  -3.98%  '

# Row 22
Set-TextValue 22 4 '6.24'
$ws.Cells.Item(22, 5).Value = '  -4.04%  '

# Row 23
$ws.Cells.Item(23, 5).Value = '  +0.39%  '

# Row 24
Set-TextValue 24 4 '65.17'
$ws.Cells.Item(24, 5).Value = '  -0.49%  '

# Row 25
$ws.Cells.Item(25, 5).Value = '  +0.46%  '

# Row 26
Set-TextValue 26 4 '0.369'
$ws.Cells.Item(26, 5).Value = '  -4.20%  '

# Row 27
Set-TextValue 27 4 '0.147'
$ws.Cells.Item(27, 5).Value = '  -4.67%  '

# Row 28
Set-TextValue 28 4 '7.21'
$ws.Cells.Item(28, 5).Value = '  -5.33%  '

# Row 29
Set-TextValue 29 4 '172.07'
$ws.Cells.Item(29, 5).Value = '  -0.86%  '

# Row 30
Set-TextValue 30 4 '0.0₃0713'
$ws.Cells.Item(30, 5).Value = '  -4.07%  '

# Row 31
Set-TextValue 31 4 '1.63'
$ws.Cells.Item(31, 5).Value = '  -4.40%  '

# Row 32
$ws.Cells.Item(32, 5).Value = '  +0.08%  '

# Row 33
$ws.Cells.Item(33, 2).Value = 'Aptos'
$ws.Cells.Item(33, 3).Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
Set-TextValue 33 4 '5.75'
$ws.Cells.Item(33, 5).Value = '  -7.94%  '

# Row 34
$ws.Cells.Item(34, 2).Value = 'FirstDigitalUSD'
$ws.Cells.Item(34, 3).Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
Set-TextValue 34 4 '0.997'
$ws.Cells.Item(34, 5).Value = '  +0.82%  '

# Row 35
$ws.Cells.Item(35, 2).Value = 'Fetch.AI'
$ws.Cells.Item(35, 3).Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue 35 4 '1.08'
$ws.Cells.Item(35, 5).Value = '  -7.85%  '

# Row 36
Set-TextValue 36 4 '17.62'
$ws.Cells.Item(36, 5).Value = '  -3.05%  '

# Row 37
$ws.Cells.Item(37, 5).Value = '  -8.42%  '

# Row 38
Set-TextValue 38 4 '3.77'
$ws.Cells.Item(38, 5).Value = '  -2.93%  '

# Row 39
Set-TextValue 39 4 '35.95'
$ws.Cells.Item(39, 5).Value = '  -1.57%  '

# Row 40
Set-TextValue 40 4 '0.793'
$ws.Cells.Item(40, 5).Value = '  -3.28%  '

# Row 41
$ws.Cells.Item(41, 5).Value = '  -4.61%  '

# Row 42
Set-TextValue 42 4 '131.19'
$ws.Cells.Item(42, 5).Value = '  -4.18%  '

# Row 43
Set-TextValue 43 4 '3.34'
$ws.Cells.Item(43, 5).Value = '  -2.83%  '

# Row 44
Set-TextValue 44 4 '4.75'
$ws.Cells.Item(44, 5).Value = '  -6.88%  '

# Row 45
Set-TextValue 45 4 '0.564'
$ws.Cells.Item(45, 5).Value = '  -1.94%  '

# Row 46
Set-TextValue 46 4 '0.0903'
$ws.Cells.Item(46, 5).Value = '  -1.44%  '

# Row 47
Set-TextValue 47 4 '241.19'
$ws.Cells.Item(47, 5).Value = '  -8.61%  '

# Row 48
Set-TextValue 48 4 '0.0479'
$ws.Cells.Item(48, 5).Value = '  -4.33%  '

# Row 49
$ws.Cells.Item(49, 5).Value = '  -4.45%  '

# Row 50
Set-TextValue 50 4 '16.95'
$ws.Cells.Item(50, 5).Value = '  -1.74%  '

# Row 51
$ws.Cells.Item(51, 5).Value = '  -3.95%  '
